$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Rename the authorization-definition field from restrict_overrides to
# extra_field_entitlements.
$ws.Range("B11").Value = "extra_field_entitlements"

# The "model" sheet tab becomes the active / selected sheet (was "survey"),
# with the cursor left on F12.
$ws.Activate()
$ws.Range("F12").Select()
